# Update cryptocurrency price/volume figures as reflected in the
# latest GitHub Actions scrape (commit: "Updated symbol list").
#
# All affected cells live in columns D (Price) and E (Volume(1h))
# and are stored as text in the workbook, so we force a text
# number format before writing the values (otherwise Excel COM
# would coerce strings like "261.73" or "1.01%" into numeric /
# percentage values), then restore the original style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "261.73"
    "E2" = "1.01%"
    "D3" = "27.19"
    "E3" = "0.52%"
    "E4" = "0.05%"
    "E5" = "3.13%"
    "D6" = "6.718"
    "D7" = "0.8500"
    "E7" = "-1.04%"
    "D8" = "0.9148"
    "E8" = "-1.36%"
    "E9" = "1.15%"
    "D10" = "0.04668"
    "E10" = "-5.56%"
    "D11" = "0.07093"
    "E11" = "0.26%"
    "D12" = "0.03151"
    "E12" = "1.83%"
    "D13" = "0.09055"
    "E13" = "-0.82%"
    "D14" = "0.001533"
    "E14" = "0.48%"
    "D15" = "0.0006167"
    "E15" = "2.00%"
    "D16" = "0.006116"
    "E16" = "0.22%"
    "D17" = "3.471"
    "E17" = "0.17%"
    "D18" = "3.171"
    "E18" = "0.05%"
    "E19" = "0.65%"
    "E21" = "1.04%"
    "D22" = "4.082"
    "E22" = "-1.37%"
    "D23" = "0.04245"
    "E23" = "-0.01%"
    "D24" = "0.001212"
    "E24" = "-0.34%"
    "D25" = "0.004133"
    "E25" = "2.39%"
    "E27" = "5.07%"
    "D40" = "0.03908"
    "E40" = "1.67%"
    "E41" = "-0.16%"
    "D42" = "0.004130"
    "E42" = "3.09%"
    "E43" = "-0.72%"
    "D45" = "0.00005172"
    "E45" = "0.79%"
    "D46" = "0.00000000751"
    "E46" = "0.08%"
    "E47" = "-34.19%"
    "D48" = "0.1667"
    "E48" = "26.16%"
    "D49" = "0.00002102"
    "E49" = "0.07%"
    "E50" = "0.08%"
}

$targetRange = $ws.Range("D2:E50")
$originalStyle = $targetRange.Style
$targetRange.NumberFormat = "@"

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$targetRange.Style = $originalStyle

Write-Output ("Updated {0} cells." -f $updates.Count)
